$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "34.607.79"
$ws.Range("E2").Value = "  +0.33%  "
Set-TextValue $ws.Range("D3") "1.812.42"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue $ws.Range("D5") "225.84"
$ws.Range("E5").Value = "  -1.25%  "
Set-TextValue $ws.Range("D6") "0.600"
$ws.Range("E6").Value = "  +3.21%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  -0.16%  "
Set-TextValue $ws.Range("D8") "37.20"
$ws.Range("E8").Value = "  +6.65%  "
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  +1.42%  "
Set-TextValue $ws.Range("D12") "2.075.49"
$ws.Range("E12").Value = "  +0.53%  "
Set-TextValue $ws.Range("D13") "11.33"
$ws.Range("E13").Value = "  +1.50%  "
Set-TextValue $ws.Range("D14") "1.817.23"
$ws.Range("E14").Value = "  +0.80%  "
Set-TextValue $ws.Range("D15") "0.635"
$ws.Range("E15").Value = "  -1.45%  "
Set-TextValue $ws.Range("D16") "34.558.94"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("E17").Value = "  +1.30%  "
Set-TextValue $ws.Range("D18") "68.76"
$ws.Range("E18").Value = "  -0.49%  "
Set-TextValue $ws.Range("D19") "243.19"
$ws.Range("E19").Value = "  -0.94%  "
Set-TextValue $ws.Range("D20") "0.0₃0778"
$ws.Range("E20").Value = "  -2.51%  "
Set-TextValue $ws.Range("D21") "11.27"
$ws.Range("E21").Value = "  -1.83%  "
Set-TextValue $ws.Range("D22") "1.00"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("E24").Value = "  +4.24%  "
Set-TextValue $ws.Range("D25") "172.41"
$ws.Range("E25").Value = "  -0.44%  "
Set-TextValue $ws.Range("D26") "7.89"
$ws.Range("E26").Value = "  -0.01%  "
Set-TextValue $ws.Range("D27") "17.28"
$ws.Range("E27").Value = "  +2.68%  "
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  -0.26%  "
Set-TextValue $ws.Range("D31") "3.94"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("E32").Value = "  -1.26%  "
Set-TextValue $ws.Range("D33") "0.0518"
$ws.Range("E33").Value = "  -2.63%  "
Set-TextValue $ws.Range("D34") "1.83"
$ws.Range("E34").Value = "  -0.61%  "
Set-TextValue $ws.Range("D35") "1.366.47"
$ws.Range("E35").Value = "  -2.14%  "
Set-TextValue $ws.Range("D36") "0.655"
$ws.Range("E36").Value = "  -4.45%  "
Set-TextValue $ws.Range("D37") "1.07"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("E38").Value = "  -4.55%  "
$ws.Range("E39").Value = "  -1.41%  "
Set-TextValue $ws.Range("D40") "2.43"
$ws.Range("E40").Value = "  +1.30%  "
Set-TextValue $ws.Range("D41") "81.18"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("E44").Value = "  +5.20%  "
$ws.Range("E45").Value = "  +1.14%  "
Set-TextValue $ws.Range("D46") "0.0501"
$ws.Range("E46").Value = "  -2.08%  "
Set-TextValue $ws.Range("D47") "1.974.72"
$ws.Range("E47").Value = "  +0.51%  "
Set-TextValue $ws.Range("D48") "5.78"
$ws.Range("E48").Value = "  -3.76%  "
Set-TextValue $ws.Range("D49") "1.00"
$ws.Range("E49").Value = "  -0.15%  "
Set-TextValue $ws.Range("D50") "102.66"
$ws.Range("E50").Value = "  -2.12%  "
Set-TextValue $ws.Range("D51") "0.0₆0121"
$ws.Range("E51").Value = "  -7.28%  "
